# Auto-generated edit script: refreshes market-board derived value columns
# (H..N) across several worksheets, matching the scheduled runner's output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1186.9487
$ws.Range("J40").Value = 1693.2
$ws.Range("L40").Value = 1693.2
$ws.Range("N40").Value = -2043.2
$ws.Range("H116").Value = 4131.8887
$ws.Range("I116").Value = 4171.5713
$ws.Range("J116").Value = 3993
$ws.Range("K116").Value = 4171.5713
$ws.Range("L116").Value = 3993
$ws.Range("M116").Value = -729.5712999999996
$ws.Range("N116").Value = -10877
$ws.Range("H137").Value = 3229.7144
$ws.Range("I137").Value = 4710.3076
$ws.Range("J137").Value = 1946.5333
$ws.Range("K137").Value = 14130.9228
$ws.Range("L137").Value = 5839.5999
$ws.Range("M137").Value = -11580.9228
$ws.Range("N137").Value = -10939.5999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1199.8
$ws.Range("I2").Value = 1199.8
$ws.Range("K2").Value = 1199.8
$ws.Range("M2").Value = -1086.8
$ws.Range("H116").Value = 1199.8
$ws.Range("I116").Value = 1199.8
$ws.Range("K116").Value = 1199.8
$ws.Range("M116").Value = 1094.2
$ws.Range("H122").Value = 3195.8
$ws.Range("I122").Value = 2925.875
$ws.Range("J122").Value = 4275.5
$ws.Range("K122").Value = 8777.625
$ws.Range("L122").Value = 12826.5
$ws.Range("M122").Value = -6327.625
$ws.Range("N122").Value = -17726.5
$ws.Range("H132").Value = 2842.7886
$ws.Range("I132").Value = 2729.4062
$ws.Range("J132").Value = 3024.2
$ws.Range("K132").Value = 8188.2186
$ws.Range("L132").Value = 9072.599999999999
$ws.Range("M132").Value = -5658.2186
$ws.Range("N132").Value = -14132.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1199.8
$ws.Range("I3").Value = 1199.8
$ws.Range("K3").Value = 1199.8
$ws.Range("M3").Value = -1085.8
$ws.Range("H86").Value = 9468.5
$ws.Range("I86").Value = 13425
$ws.Range("J86").Value = 1555.5
$ws.Range("K86").Value = 13425
$ws.Range("L86").Value = 1555.5
$ws.Range("M86").Value = -12302
$ws.Range("N86").Value = -3801.5
$ws.Range("H89").Value = 9468.5
$ws.Range("I89").Value = 13425
$ws.Range("J89").Value = 1555.5
$ws.Range("K89").Value = 67125
$ws.Range("L89").Value = 7777.5
$ws.Range("M89").Value = -61509
$ws.Range("N89").Value = -19009.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 35950
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 35950
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 35950
$ws.Range("N68").Value = -37448
$ws.Range("M68").Value = $null
$ws.Range("H70").Value = 8985.714
$ws.Range("J70").Value = 8985.714
$ws.Range("L70").Value = 8985.714
$ws.Range("N70").Value = -9615.714
$ws.Range("H71").Value = 35950
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 35950
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 107850
$ws.Range("N71").Value = -115338
$ws.Range("M71").Value = $null
$ws.Range("H73").Value = 8985.714
$ws.Range("J73").Value = 8985.714
$ws.Range("L73").Value = 8985.714
$ws.Range("N73").Value = -11169.714
$ws.Range("H74").Value = 23771.334
$ws.Range("J74").Value = 23771.334
$ws.Range("L74").Value = 23771.334
$ws.Range("N74").Value = -25519.334
$ws.Range("H77").Value = 23771.334
$ws.Range("J77").Value = 23771.334
$ws.Range("L77").Value = 71314.00199999999
$ws.Range("N77").Value = -80050.00199999999
$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -42246
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -131232
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
$ws.Range("H122").Value = 1130.8334
$ws.Range("I122").Value = 959.3333
$ws.Range("J122").Value = 1302.3334
$ws.Range("K122").Value = 2877.9999
$ws.Range("L122").Value = 3907.0002
$ws.Range("M122").Value = -427.9998999999998
$ws.Range("N122").Value = -8807.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4496.0938
$ws.Range("I70").Value = 4198.4546
$ws.Range("J70").Value = 4652
$ws.Range("K70").Value = 4198.4546
$ws.Range("L70").Value = 4652
$ws.Range("M70").Value = -3928.4546
$ws.Range("N70").Value = -5192
$ws.Range("H73").Value = 4496.0938
$ws.Range("I73").Value = 4198.4546
$ws.Range("J73").Value = 4652
$ws.Range("K73").Value = 4198.4546
$ws.Range("L73").Value = 4652
$ws.Range("M73").Value = -3262.4546
$ws.Range("N73").Value = -6524
$ws.Range("H113").Value = 2603.4
$ws.Range("I113").Value = 2416.8333
$ws.Range("J113").Value = 2883.25
$ws.Range("K113").Value = 2416.8333
$ws.Range("L113").Value = 2883.25
$ws.Range("M113").Value = -246.8332999999998
$ws.Range("N113").Value = -7223.25
$ws.Range("H122").Value = 1633.375
$ws.Range("I122").Value = 1275
$ws.Range("J122").Value = 1991.75
$ws.Range("K122").Value = 3825
$ws.Range("L122").Value = 5975.25
$ws.Range("M122").Value = -1375
$ws.Range("N122").Value = -10875.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 750
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6840
$ws.Range("H122").Value = 2041.7333
$ws.Range("I122").Value = 2081.3333
$ws.Range("J122").Value = 1883.3334
$ws.Range("K122").Value = 6243.999899999999
$ws.Range("L122").Value = 5650.0002
$ws.Range("M122").Value = -3793.999899999999
$ws.Range("N122").Value = -10550.0002
$ws.Range("H136").Value = 13348721
$ws.Range("I136").Value = 26343804
$ws.Range("J136").Value = 2419.5945
$ws.Range("K136").Value = 79031412
$ws.Range("L136").Value = 7258.7835
$ws.Range("M136").Value = -79028862
$ws.Range("N136").Value = -12358.7835
